$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1) Update odds for the existing fixture in row 2 (Newells vs Independiente)
# ---------------------------------------------------------------
$ws.Range("G2").Value = 3.5
$ws.Range("I2").Value = 2.5
$ws.Range("J2").Value = 4.5
$ws.Range("L2").Value = 3.5
$ws.Range("M2").Value = 1.18
$ws.Range("N2").Value = 4.5
$ws.Range("Q2").Value = 3.6
$ws.Range("R2").Value = 1.29
$ws.Range("W2").Value = 6.5
$ws.Range("X2").Value = 15
$ws.Range("Y2").Value = 15
$ws.Range("AB2").Value = 51
$ws.Range("AD2").Value = 5.5
$ws.Range("AH2").Value = 10
$ws.Range("AJ2").Value = 26
$ws.Range("AO2").Value = 23
$ws.Range("AQ2").Value = 81
$ws.Range("AZ2").Value = 67

# ---------------------------------------------------------------
# 2) Insert a brand-new fixture row at row 3 (QRwclaLb: Bolivia - Guabira vs Tomayapo)
#    This shifts the former rows 3,4,5 down to 4,5,6.
# ---------------------------------------------------------------
$ws.Rows("3").Insert()

$newRow3 = New-Object 'object[,]' 1,56
$newRow3[0,0] = 'QRwclaLb'
$newRow3[0,1] = '26/11/2024'
$newRow3[0,2] = '20:30'
$newRow3[0,3] = 'BOLIVIA - DIVISION PROFESIONAL'
$newRow3[0,4] = 'Guabira'
$newRow3[0,5] = 'Tomayapo'
$newRow3[0,6] = 1.75
$newRow3[0,7] = 3.8
$newRow3[0,8] = 4.33
$newRow3[0,9] = 2.38
$newRow3[0,10] = 2.3
$newRow3[0,11] = 4.5
$newRow3[0,12] = 1.04
$newRow3[0,13] = 13
$newRow3[0,14] = 1.22
$newRow3[0,15] = 4
$newRow3[0,16] = 1.75
$newRow3[0,17] = 2.05
$newRow3[0,18] = 1.3
$newRow3[0,19] = 3.25
$newRow3[0,20] = 1.73
$newRow3[0,21] = 2
$newRow3[0,22] = 8
$newRow3[0,23] = 9
$newRow3[0,24] = 8.5
$newRow3[0,25] = 15
$newRow3[0,26] = 13
$newRow3[0,27] = 23
$newRow3[0,28] = 12
$newRow3[0,29] = 7
$newRow3[0,30] = 15
$newRow3[0,31] = 41
$newRow3[0,32] = 13
$newRow3[0,33] = 23
$newRow3[0,34] = 15
$newRow3[0,35] = 41
$newRow3[0,36] = 34
$newRow3[0,37] = 41
$newRow3[0,38] = 201
$newRow3[0,39] = 3.75
$newRow3[0,40] = 9
$newRow3[0,41] = 19
$newRow3[0,42] = 29
$newRow3[0,43] = 51
$newRow3[0,44] = 126
$newRow3[0,45] = 3.25
$newRow3[0,46] = 8
$newRow3[0,47] = 51
$newRow3[0,48] = 6
$newRow3[0,49] = 23
$newRow3[0,50] = 29
$newRow3[0,51] = 81
$newRow3[0,52] = 81
$newRow3[0,53] = 201
$newRow3[0,54] = ""
$newRow3[0,55] = ""
$ws.Range("A3:BD3").Value = $newRow3

# ---------------------------------------------------------------
# 3) Update odds for the fixture that shifted from row 3 -> row 4 (dfKzglAQ: Fluminense vs Criciuma)
# ---------------------------------------------------------------
$ws.Range("M4").Value = 1.08
$ws.Range("N4").Value = 8
$ws.Range("Q4").Value = 2.15
$ws.Range("R4").Value = 1.67
$ws.Range("S4").Value = 1.41
$ws.Range("T4").Value = 2.62
$ws.Range("AC4").Value = 8
$ws.Range("AM4").Value = 451

# ---------------------------------------------------------------
# 4) Update odds for the fixture that shifted from row 4 -> row 5 (8nJEo620: Fortaleza vs Flamengo RJ)
# ---------------------------------------------------------------
$ws.Range("N5").Value = 9

# ---------------------------------------------------------------
# 5) Update odds for the fixture that shifted from row 5 -> row 6 (UeSMa6cR: Wanderers vs Miramar)
# ---------------------------------------------------------------
$ws.Range("G6").Value = 2.55
$ws.Range("H6").Value = 3.3
$ws.Range("I6").Value = 2.7
$ws.Range("O6").Value = 1.3
$ws.Range("P6").Value = 3.4
$ws.Range("Q6").Value = 2.03
$ws.Range("R6").Value = 1.83
$ws.Range("X6").Value = 13
$ws.Range("Y6").Value = 10
$ws.Range("Z6").Value = 23
$ws.Range("AA6").Value = 21
$ws.Range("AG6").Value = 9
$ws.Range("AH6").Value = 13
$ws.Range("AI6").Value = 10
$ws.Range("AJ6").Value = 26
$ws.Range("AK6").Value = 21
$ws.Range("AL6").Value = 29
$ws.Range("AN6").Value = 4.5
